$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (column C) date from 45442 to 45443 for all existing
# data rows (rows 2-28).
for ($r = 2; $r -le 28; $r++) {
    $ws.Cells.Item($r, 3).Value = 45443
}

# Row 28 picks up an explicit row height (15, custom) as part of the edit.
$ws.Rows.Item(28).RowHeight = 15

# Append the new record as row 29.
$ws.Cells.Item(29, 1).Value = "A 21829-2024"

$ws.Cells.Item(29, 2).Value = 45442
$ws.Cells.Item(29, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(29, 3).Value = 45443
$ws.Cells.Item(29, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(29, 4).Value = "OKÄNT"
$ws.Cells.Item(29, 5).Value = "OKÄNT"

$ws.Cells.Item(29, 7).Value = 1.1
$ws.Cells.Item(29, 8).Value = 0
$ws.Cells.Item(29, 9).Value = 0
$ws.Cells.Item(29, 10).Value = 0
$ws.Cells.Item(29, 11).Value = 0
$ws.Cells.Item(29, 12).Value = 0
$ws.Cells.Item(29, 13).Value = 0
$ws.Cells.Item(29, 14).Value = 0
$ws.Cells.Item(29, 15).Value = 0
$ws.Cells.Item(29, 16).Value = 0
$ws.Cells.Item(29, 17).Value = 0

# Column R keeps the blank-but-wrap-formatted style used throughout the sheet.
$ws.Cells.Item(29, 18).WrapText = $true
